# Insert a new data row at row 31 (pushing the existing rows 31-77 down to 32-78)
# and populate it with a new Maracuyá price record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("31:31").Insert()

$ws.Cells.Item(31, 1).Value  = 10
$ws.Cells.Item(31, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(31, 3).Value  = "La Araucanía"
$ws.Cells.Item(31, 4).Value  = 44994
$ws.Cells.Item(31, 5).Value  = 9
$ws.Cells.Item(31, 6).Value  = "Fruta"
$ws.Cells.Item(31, 7).Value  = 100108
$ws.Cells.Item(31, 8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(31, 9).Value  = 100108003
$ws.Cells.Item(31, 10).Value = "Maracuyá"
$ws.Cells.Item(31, 11).Value = "Sin especificar"
$ws.Cells.Item(31, 12).Value = "Primera"
$ws.Cells.Item(31, 13).Value = 30
$ws.Cells.Item(31, 14).Value = 60000
$ws.Cells.Item(31, 15).Value = 60000
$ws.Cells.Item(31, 16).Value = 60000
$ws.Cells.Item(31, 17).Value = "$/caja 18 kilos"
$ws.Cells.Item(31, 18).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(31, 19).Value = 3333
$ws.Cells.Item(31, 20).Value = 18
